{"js": "// Replace the date line and each of the 26 division-exercise answer\n// cells with their updated values, per the commit diff. Every old\n// string in the document occurs exactly once, so a simple\n// search-and-replace per pair is sufficient and keeps all paragraph /\n// run formatting (fonts, size, etc.) untouched.\nconst replacements = [\n  [\"2025-06-25 Wednesday\", \"2025-06-26 Thursday\"],\n  [\"471\u00f74=117, 3\", \"768\u00f73=256, 0\"],\n  [\"681\u00f78=85, 1\", \"915\u00f74=228, 3\"],\n  [\"898\u00f77=128, 2\", \"442\u00f77=63, 1\"],\n  [\"891\u00f76=148, 3\", \"974\u00f74=243, 2\"],\n  [\"677\u00f75=135, 2\", \"982\u00f75=196, 2\"],\n  [\"698\u00f79=77, 5\", \"584\u00f75=116, 4\"],\n  [\"382\u00f76=63, 4\", \"791\u00f72=395, 1\"],\n  [\"938\u00f72=469, 0\", \"911\u00f73=303, 2\"],\n  [\"739\u00f77=105, 4\", \"831\u00f74=207, 3\"],\n  [\"124\u00f78=15, 4\", \"162\u00f72=81, 0\"],\n  [\"354\u00f78=44, 2\", \"233\u00f72=116, 1\"],\n  [\"307\u00f75=61, 2\", \"846\u00f72=423, 0\"],\n  [\"351\u00f75=70, 1\", \"477\u00f73=159, 0\"],\n  [\"944\u00f75=188, 4\", \"561\u00f75=112, 1\"],\n  [\"789\u00f74=197, 1\", \"722\u00f73=240, 2\"],\n  [\"395\u00f74=98, 3\", \"686\u00f78=85, 6\"],\n  [\"818\u00f74=204, 2\", \"221\u00f75=44, 1\"],\n  [\"696\u00f72=348, 0\", \"888\u00f76=148, 0\"],\n  [\"327\u00f75=65, 2\", \"291\u00f74=72, 3\"],\n  [\"448\u00f79=49, 7\", \"751\u00f75=150, 1\"],\n  [\"758\u00f75=151, 3\", \"236\u00f77=33, 5\"],\n  [\"800\u00f73=266, 2\", \"971\u00f76=161, 5\"],\n  [\"288\u00f72=144, 0\", \"770\u00f78=96, 2\"],\n  [\"760\u00f72=380, 0\", \"611\u00f72=305, 1\"],\n  [\"387\u00f78=48, 3\", \"250\u00f75=50, 0\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each of the 26 division-exercise answer\n# cells with their updated values, per the commit diff. Every old\n# string in the document occurs exactly once, so a Find/Replace pass\n# per pair is sufficient and preserves existing run/paragraph\n# formatting (fonts, size, etc.).\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @{old=\"2025-06-25 Wednesday\"; new=\"2025-06-26 Thursday\"},\n  @{old=\"471\u00f74=117, 3\"; new=\"768\u00f73=256, 0\"},\n  @{old=\"681\u00f78=85, 1\"; new=\"915\u00f74=228, 3\"},\n  @{old=\"898\u00f77=128, 2\"; new=\"442\u00f77=63, 1\"},\n  @{old=\"891\u00f76=148, 3\"; new=\"974\u00f74=243, 2\"},\n  @{old=\"677\u00f75=135, 2\"; new=\"982\u00f75=196, 2\"},\n  @{old=\"698\u00f79=77, 5\"; new=\"584\u00f75=116, 4\"},\n  @{old=\"382\u00f76=63, 4\"; new=\"791\u00f72=395, 1\"},\n  @{old=\"938\u00f72=469, 0\"; new=\"911\u00f73=303, 2\"},\n  @{old=\"739\u00f77=105, 4\"; new=\"831\u00f74=207, 3\"},\n  @{old=\"124\u00f78=15, 4\"; new=\"162\u00f72=81, 0\"},\n  @{old=\"354\u00f78=44, 2\"; new=\"233\u00f72=116, 1\"},\n  @{old=\"307\u00f75=61, 2\"; new=\"846\u00f72=423, 0\"},\n  @{old=\"351\u00f75=70, 1\"; new=\"477\u00f73=159, 0\"},\n  @{old=\"944\u00f75=188, 4\"; new=\"561\u00f75=112, 1\"},\n  @{old=\"789\u00f74=197, 1\"; new=\"722\u00f73=240, 2\"},\n  @{old=\"395\u00f74=98, 3\"; new=\"686\u00f78=85, 6\"},\n  @{old=\"818\u00f74=204, 2\"; new=\"221\u00f75=44, 1\"},\n  @{old=\"696\u00f72=348, 0\"; new=\"888\u00f76=148, 0\"},\n  @{old=\"327\u00f75=65, 2\"; new=\"291\u00f74=72, 3\"},\n  @{old=\"448\u00f79=49, 7\"; new=\"751\u00f75=150, 1\"},\n  @{old=\"758\u00f75=151, 3\"; new=\"236\u00f77=33, 5\"},\n  @{old=\"800\u00f73=266, 2\"; new=\"971\u00f76=161, 5\"},\n  @{old=\"288\u00f72=144, 0\"; new=\"770\u00f78=96, 2\"},\n  @{old=\"760\u00f72=380, 0\"; new=\"611\u00f72=305, 1\"},\n  @{old=\"387\u00f78=48, 3\"; new=\"250\u00f75=50, 0\"}\n)\n\nforeach ($p in $pairs) {\n  $range = $d.Content\n  $range.Find.ClearFormatting()\n  $range.Find.Replacement.ClearFormatting()\n  $range.Find.Text = $p.old\n  $range.Find.Replacement.Text = $p.new\n  $range.Find.Execute($p.old, $false, $false, $false, $false, $false, $true, 1, $false, $p.new, 2)\n}\n"}
